# Sheet3 ("存款" / deposits): add bank/deposit_type/currency detail columns
# G:M (property_category, category, date, legislator_name, legislator_id,
# source_file, index) mirroring the pattern already used on the other
# sheets (e.g. "股票"). Columns A:F keep their existing values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Header row (row 1) ----
$ws.Cells.Item(1,7).Value  = "property_category"
$ws.Cells.Item(1,8).Value  = "category"
$ws.Cells.Item(1,9).Value  = "date"
$ws.Cells.Item(1,10).Value = "legislator_name"
$ws.Cells.Item(1,11).Value = "legislator_id"
$ws.Cells.Item(1,12).Value = "source_file"
$ws.Cells.Item(1,13).Value = "index"

# Match the bold/boxed header style used by B1:F1 (style copied from B1).
$ws.Range("B1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("M1").PasteSpecial(-4122)

# ---- Data rows (rows 2-7), one deposit record per row ----
# Maps worksheet row -> the "index" value (same as column A on that row).
$indexByRow = @{ 2 = 78; 3 = 79; 4 = 80; 5 = 81; 6 = 82; 7 = 83 }

foreach ($r in 2..7) {
    $ws.Cells.Item($r, 7).Value  = "deposit"
    $ws.Cells.Item($r, 8).Value  = "normal"
    # Leading apostrophe forces text (otherwise "2012-04-18" gets parsed as
    # a date serial); re-copying the plain data-cell format afterwards
    # clears the quote-prefix marker Excel adds for that.
    $ws.Cells.Item($r, 9).Value  = "'2012-04-18"
    $ws.Cells.Item($r, 10).Value = "林明溱"
    $ws.Cells.Item($r, 11).Value = 1706
    $ws.Cells.Item($r, 12).Value = "tmp80511"
    $ws.Cells.Item($r, 13).Value = $indexByRow[$r]
}

$ws.Range("B2").Copy()
foreach ($r in 2..7) {
    $ws.Cells.Item($r, 9).PasteSpecial(-4122)
}
